# class & lambda slides
#
# The title on the "Some methods many class have" slide (SlideID 619) was
# split across two runs:
#   run 1: "Some methods "
#   run 2: "many class have"
# The author retyped it as a single corrected sentence:
#   "Some methods many classes have"
# merging it back down into one run.

$p = $ppt.ActivePresentation

$target = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    if ($slide.SlideID -eq 619) {
        $target = $slide
        break
    }
}
if ($target -eq $null) {
    throw "Could not find slide with SlideID 619"
}

$title = $target.Shapes.Item(1)
$tr = $title.TextFrame.TextRange

$oldText = "Some methods many class have"
$newText = "Some methods many classes have"
$firstRunText = "Some methods "

if ($tr.Text -eq $oldText) {
    # Drop the leading "Some methods " run so only the second run (which is
    # already marked dirty="0") survives, then restore the full corrected
    # sentence in one go -- this collapses the paragraph back down to a
    # single run, matching how PowerPoint merges runs when you retype a
    # whole line rather than editing it in place.
    $lead = $tr.Characters(1, $firstRunText.Length)
    $lead.Text = ""
    $tr.Text = $newText
} else {
    # Fallback in case the text has already been touched.
    $tr.Text = $newText
}
